# "found/fixed typo in resume"
#
# Original:
#   Updated webserver (JavaScript) for remote access to optical
#   polarization ray tracing software
#
# New:
#   Updated webserver (node.js) for remote access to optical
#   polarization ray tracing software
#
# In addition, Word's auto-managed "_GoBack" bookmark (which always
# marks the location of the most recent edit) moves from the end of
# the ", Polaris-M" run (its old location) to right after the newly
# typed "node.js" text.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the word "JavaScript" inside the unique anchor phrase
# "Updated webserver (JavaScript)" so we don't touch the unrelated
# "JavaScript" mention later in the Programming Languages list.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Updated webserver (JavaScript)", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
if (-not $anchor.Find.Found) {
    throw "Could not find the 'Updated webserver (JavaScript)' phrase"
}

# Range covering just the word "JavaScript" within that phrase.
$jsStart = $anchor.Start + "Updated webserver (".Length
$jsEnd = $jsStart + "JavaScript".Length
$jsRange = $d.Range($jsStart, $jsEnd)
if ($jsRange.Text -ne "JavaScript") {
    throw "Unexpected text at computed JavaScript range: '$($jsRange.Text)'"
}

# ------------------------------------------------------------------
# Drop a throwaway bookmark right before "JavaScript" first - this
# keeps "Updated webserver (" as its own run once the replacement
# text goes in right after it, instead of the two being silently
# coalesced back into a single run.
# ------------------------------------------------------------------
$splitPoint = $d.Range($jsStart, $jsStart)
$d.Bookmarks.Add("zzTempSplit", $splitPoint)

# Replace "JavaScript" with "node.js".
$jsRange = $d.Range($jsStart, $jsEnd)
$jsRange.Text = "node.js"

# Re-typing the new word over a fresh range also makes Word drop the
# (inherited) xml:space="preserve" flag it has no reason to carry.
$nodeRange = $d.Range($jsStart, $jsStart + "node.js".Length)
$nodeRange.Text = ""
$reinsert = $d.Range($jsStart, $jsStart)
$reinsert.InsertBefore("node.js")

# ------------------------------------------------------------------
# Move the "_GoBack" bookmark (Word always keeps exactly one, at the
# site of the latest edit) to just after "node.js".
# ------------------------------------------------------------------
$goBackPos = $jsStart + "node.js".Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Drop the temporary split-marker bookmark now that "_GoBack" is the
# only thing still holding the run boundary open.
$d.Bookmarks.Item("zzTempSplit").Delete()

# ------------------------------------------------------------------
# Finally, re-insert "Updated webserver (" itself so it likewise
# loses any inherited xml:space="preserve" it doesn't need.
# ------------------------------------------------------------------
$prefixStart = $anchor.Start
$prefixText = "Updated webserver ("
$prefixRange = $d.Range($prefixStart, $prefixStart + $prefixText.Length)
if ($prefixRange.Text -ne $prefixText) {
    throw "Unexpected prefix text: '$($prefixRange.Text)'"
}
$prefixRange.Text = ""
$prefixReinsert = $d.Range($prefixStart, $prefixStart)
$prefixReinsert.InsertBefore($prefixText)

Write-Host "Replaced JavaScript -> node.js and relocated _GoBack bookmark"
